$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.494.31"
$ws.Range("E2").Value = "  -1.29%  "

$ws.Range("D3").Value = "2.485.74"
$ws.Range("E3").Value = "  -1.56%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.92"
$ws.Range("E5").Value = "  -2.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.04"
$ws.Range("E6").Value = "  -3.22%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("E9").Value = "  -1.00%  "

$ws.Range("E10").Value = "  -2.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.40"
$ws.Range("E11").Value = "  +0.52%  "

$ws.Range("E12").Value = "  -0.96%  "

$ws.Range("D13").Value = "2.928.48"
$ws.Range("E13").Value = "  -1.59%  "

$ws.Range("D14").Value = "58.407.69"
$ws.Range("E14").Value = "  -1.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.47"

$ws.Range("E16").Value = "  -1.97%  "

$ws.Range("D17").Value = "2.488.47"
$ws.Range("E17").Value = "  -1.37%  "

$ws.Range("E18").Value = "  -1.56%  "

$ws.Range("E19").Value = "  -1.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "321.83"
$ws.Range("E20").Value = "  -1.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.84"
$ws.Range("E22").Value = "  -1.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.37"
$ws.Range("E23").Value = "  -2.07%  "

$ws.Range("E24").Value = "  -2.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("E26").Value = "  -2.53%  "

$ws.Range("E27").Value = "  -2.49%  "

$ws.Range("E28").Value = "  -3.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.42"
$ws.Range("E29").Value = "  -5.22%  "

$ws.Range("E30").Value = "  -3.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.11"
$ws.Range("E31").Value = "  -2.74%  "

$ws.Range("E32").Value = "  -4.98%  "

$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.28"
$ws.Range("E35").Value = "  -1.41%  "

$ws.Range("E36").Value = "  -8.61%  "

$ws.Range("E37").Value = "  -3.03%  "

$ws.Range("E38").Value = "  -4.05%  "

$ws.Range("E39").Value = "  -2.63%  "

$ws.Range("E40").Value = "  -3.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "276.04"
$ws.Range("E41").Value = "  -2.95%  "

$ws.Range("E42").Value = "  -5.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.597"
$ws.Range("E43").Value = "  -1.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "127.48"
$ws.Range("E44").Value = "  -2.45%  "

$ws.Range("E45").Value = "  -1.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0496"
$ws.Range("E46").Value = "  -2.99%  "

$ws.Range("E47").Value = "  -2.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "17.25"
$ws.Range("E48").Value = "  -1.59%  "

$ws.Range("D49").Value = "1.741.18"
$ws.Range("E49").Value = "  -1.43%  "

$ws.Range("E50").Value = "  -1.50%  "

$ws.Range("E51").Value = "  -1.78%  "
